$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Услуга) to hold "Лицевой счет" (Account Number)
$ws.Columns("D").Insert()

# Header
$ws.Range("D1").Value = "Лицевой счет"

# Row 2 placeholder
$ws.Range("D2").Value = "{d.meter[i].accountNumber}"

# Row 3 placeholder
$ws.Range("D3").Value = "{d.meter[i + 1].accountNumber}"
